$d = $word.ActiveDocument
$r = $d.Content.Duplicate
$found = $r.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(" (")
$r.Collapse(0)
$r.InsertAfter("Changed main")
$r.Collapse(0)
$r.InsertAfter(")")
